# Insert a new weekly price record as row 328 on the single data sheet.
# All existing rows from 328 downward shift down by one (328->329, ..., 424->425),
# and the sheet's used-range dimension grows from R424 to R425 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 328..424 down to 329..425 by inserting a blank row at 328.
$ws.Rows.Item(328).Insert()

# Populate the newly inserted row 328 with the new weekly record.
$ws.Cells.Item(328, 1).Value  = 1
$ws.Cells.Item(328, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(328, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(328, 4).Value  = 44841
$ws.Cells.Item(328, 5).Value  = 15
$ws.Cells.Item(328, 6).Value  = 100112032
$ws.Cells.Item(328, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(328, 8).Value  = "Huracán"
$ws.Cells.Item(328, 9).Value  = "Segunda"
$ws.Cells.Item(328, 10).Value = 140
$ws.Cells.Item(328, 11).Value = 15000
$ws.Cells.Item(328, 12).Value = 16000
$ws.Cells.Item(328, 13).Value = 15500
$ws.Cells.Item(328, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(328, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(328, 16).Value = 155
$ws.Cells.Item(328, 17).Value = 100
$ws.Cells.Item(328, 18).Value = "Hortaliza"
